# Update the dSF column (F) values for rows 2-6 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = 0
